$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old "Weekly Total" row for the Nov13th-Nov20th week sat at row 23 with a
# blank row 22 above it. Harvey removed that blank spacer row so the total
# moves up to row 22 (formula refs auto-adjust: SUM(C19:C21) stays the same,
# and the downstream "Current Total" formula will pick up C22 instead of C23).
$ws.Rows(22).Delete()

# The old "Current Total:" row (now shifted up to row 26) gets pushed further
# down the sheet to make room for the new week's block; clear its contents
# here (formatting is left in place - it gets reused below) and rewrite it
# at its final resting place, row 34.
$ws.Range("B26:C26").ClearContents()

# New week section: "Week:Nov20th-Nov27th"
$ws.Range("A23").Value = "Week:Nov20th-Nov27th"
$ws.Range("B24").Value = "Task"

$ws.Range("B25").Value = "Prototpying JSON Generation"
$ws.Range("C25").Value = 0.08333333333333333
$ws.Range("C25").NumberFormat = "h:mm"

$ws.Range("B26").Value = "Documentation for the Design Spec"
$ws.Range("C26").Value = 0.020833333333333332

$ws.Range("B27").Value = "Meetings"
$ws.Range("C27").Value = 1

# Weekly total for the new week (rows 28-30 left blank for future entries,
# same gap pattern as the earlier weekly blocks).
$ws.Range("B31").Value = "Weekly Total"
$ws.Range("C31").Formula = "=SUM(C25:C30)"
$ws.Range("C31").NumberFormat = "h:mm"

# Current (running) total now folds in the new week's subtotal too.
$ws.Range("B34").Value = "Current Total:"
$ws.Range("C34").Formula = "=SUM(C9+C16+C22+C31)"

# Match the author's last selection when they saved the file.
[void]$ws.Range("E30").Select()
